$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.571.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -9.27%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.621.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -9.16%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.40%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -9.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.611.83"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -9.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -11.38%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.686"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -13.69%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -16.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.88"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -13.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000282"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -16.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.24"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -12.23%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.212.31"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -8.93%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.639.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.78%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.75%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.89"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -9.73%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -13.22%  "

# Row 20
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.09"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -12.24%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "66.014.27"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -9.99%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -15.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.31"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -11.42%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -11.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.95"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -13.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.29"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -14.49%  "

# Row 27
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.84%  "

# Row 28
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.93"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.62"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -14.61%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.16"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -13.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -13.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.39"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.14"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -13.38%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.52"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -9.73%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.113"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -13.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.55"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -15.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "571.05"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -10.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0871"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -17.57%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.14%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.36%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.381"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -12.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.128"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -13.48%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -11.18%  "

# Row 44
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.89"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -15.05%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0422"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -13.48%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.46"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -16.95%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -12.95%  "

# Row 49
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.39%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.644.40"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.17%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.10"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.96%  "
